$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2021-12-08"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 5.85
